# Add a new data row (row 3) to the "Site Total Alarms" sheet for 2024-02-22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date label in column A, stored as text (matches existing A2 which is text).
# Leading apostrophe forces Excel to treat the numeric-looking string as text.
$ws.Range("A3").Value = "'20240222"

# Numeric alarm counts for ARN, AZC, PKV, SDU, STL, THL columns.
$ws.Cells.Item(3, 2).Value = 103
$ws.Cells.Item(3, 3).Value = 625
$ws.Cells.Item(3, 4).Value = 345
$ws.Cells.Item(3, 5).Value = 84
$ws.Cells.Item(3, 6).Value = 45
$ws.Cells.Item(3, 7).Value = 92

# Copy the formatting from A2 (bold, bordered, centered) onto A3 so the new
# date cell matches the style of the existing date cells.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
